$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (13) recording the "14.3.14 Treffen" (meeting) entry.
# Date 41711 = 2014-03-13 (xlsx serial date), 6 hours, Paul/Kevin/Maxi/Phil all "x" (present),
# with a long comment describing next tasks.

$ws.Range("A13").Value = 41711
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("A13").HorizontalAlignment = $ws.Range("A12").HorizontalAlignment

$ws.Range("B13").Value = 6
$ws.Range("B13").HorizontalAlignment = $ws.Range("B12").HorizontalAlignment

$ws.Range("C13").Value = "x"
$ws.Range("C13").HorizontalAlignment = $ws.Range("C12").HorizontalAlignment
$ws.Range("C13").VerticalAlignment = $ws.Range("C12").VerticalAlignment

$ws.Range("D13").Value = "x"
$ws.Range("D13").HorizontalAlignment = $ws.Range("D12").HorizontalAlignment
$ws.Range("D13").VerticalAlignment = $ws.Range("D12").VerticalAlignment

$ws.Range("E13").Value = "x"
$ws.Range("E13").HorizontalAlignment = $ws.Range("E12").HorizontalAlignment
$ws.Range("E13").VerticalAlignment = $ws.Range("E12").VerticalAlignment

$ws.Range("F13").Value = "x"
$ws.Range("F13").HorizontalAlignment = $ws.Range("F12").HorizontalAlignment
$ws.Range("F13").VerticalAlignment = $ws.Range("F12").VerticalAlignment

$ws.Range("H13").Value = "login`n addfriend`n game-logic (answer questions from gameOverview, send roundResult)`ngameOverview encolourAllQuestions`ndisable ""Spielen""-button if not WaitingFor!"
$ws.Range("H13").WrapText = $true

$ws.Rows.Item(13).RowHeight = 90

# Move the view / selection down to the newly added row, like the author did.
$ws.Range("I13").Select()
